# Update the "Installation" cost category/description labels to
# "Installation and Commissioning" throughout the Inputs table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Python test sheet")

$rows = @(8, 11, 14, 23, 26, 29)

foreach ($r in $rows) {
    $descCell = $ws.Range("E" + $r)
    $descText = $descCell.Value2
    $descCell.Value = ($descText -replace "\(Installation\)$", "(Installation and Commissioning)")

    $catCell = $ws.Range("H" + $r)
    $catCell.Value = "Installation and Commissioning"
}

$ws.Range("G19").Select()
